$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 57.9
$ws.Range("N2").Value = 54.83846622768671

$ws.Range("K3").Value = 54.5
$ws.Range("N3").Value = 54.83846622768671
